$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

$ws.Range("B11").Value = 0.1220857228185415
$ws.Range("B12").Value = 0.3595817908248243
$ws.Range("C12").Value = "{'codebleu': 0.35958179082482433, 'ngram_match_score': 0.1217292643494159, 'weighted_ngram_match_score': 0.13560242383675933, 'syntax_match_score': 0.5692307692307692, 'dataflow_match_score': 0.611764705882353}"
$ws.Range("B13").Value = 0.4607057954493304
